$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.962

$ws.Range("A4").Value = -20.858
$ws.Range("B4").Value = 7.047999999999999
$ws.Range("E4").Value = 12.977

$ws.Range("B5").Value = 6.633

$ws.Range("A6").Value = -21.108

$ws.Range("A7").Value = -21.038

$ws.Range("B8").Value = 6.347

$ws.Range("E9").Value = 13.018

$ws.Range("E11").Value = 12.852

$ws.Range("E14").Value = 12.964

$ws.Range("A16").Value = -20.91
$ws.Range("B16").Value = 6.840000000000001

$ws.Range("E18").Value = 12.596

$ws.Range("A20").Value = -22.121

$ws.Range("B22").Value = 6.569

$ws.Range("E25").Value = 13.018
